# Relatorio 3 - final edits to the "Analise dos Resultados" section:
#  - merge the bookmark-split "4-PAM" sentence into a single run
#  - split the "O script extra..." paragraph in two, rewriting the tail
#    about the 4-PAM matched filter plot
#  - move the _GoBack bookmark down to the end of the document
$d = $word.ActiveDocument

# Anchor on the first sentence that changes ("O caso do 4-PAM segue...") and
# walk back one character to also capture the leading tab run.
$anchorStart = $d.Content
$null = $anchorStart.Find.Execute("O caso do 4-PAM segue", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$editStart = $anchorStart.Start - 1

# Anchor on the last sentence of the section ("...fidedignidade.") to find
# where the edited span ends.
$anchorEnd = $d.Content
$null = $anchorEnd.Find.Execute("Com valores menores, a figura tendia a ser representada com quase total fidedignidade.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$editEnd = $anchorEnd.End

$editRange = $d.Range($editStart, $editEnd)
$editRange.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>O caso do 4-PAM segue o mesmo raciocínio do sinal antipodal, com alterações no bloco de decisão e ajustes nos níveis de amplitude. Nota-se um menor desempenho em relação ao 2-PAM e ao ortogonal.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>Em comparação aos resultados das simulações analógicas, foi possível notar as mesmas características entre os sinais ortogonal, antipodal e 4-PAM. Tais comparações foram feitas qualitativamente, considerando a mesma referência da energia de bit por densidade espectral de ruído (SNR).</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>O script extra foi executado a fim de visualizar os efeitos do ruído em uma imagem</w:t></w:r><w:r><w:t xml:space="preserve"> bitmap preto e branco. Pôde-se notar que a decisão do bit, ligada diretamente à cor projetada, apresentava maiores taxas de erro com SNR menor, logo a imagem apresentava muitos pontos com falhas. </w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="708"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Os filtros casados foram implementados para transmissão em sinal ortogonal e antipodal, observando a vantagem do antipodal em relação ao outro. O filtro para 4-PAM funcionou e apresentou as características esperadas, porém por </w:t></w:r><w:r><w:t>razão</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>da manipulação das matrizes</w:t></w:r><w:r><w:t xml:space="preserve"> o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>plot</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> final da figura fic</w:t></w:r><w:r><w:t>ou</w:t></w:r><w:r><w:t xml:space="preserve"> deslocado</w:t></w:r><w:r><w:t xml:space="preserve"> e invertido.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>O impacto dos valores de BER distintos representava em aparência de “chuviscos” na imagem, dando noção visual de ruído. Com valores menores, a figura tendia a ser representada com quase total fidedignidade.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
